# Remove the decorative "separator" paragraphs (long horizontal rule made of
# U+2500 box-drawing characters, spacing before=120/after=120 twips i.e. 6pt)
# and the empty "spacer" paragraphs (no runs, spacing before=40 twips i.e. 2pt)
# that previously padded the space between tables / before section headings.
#
# We scan $d.Paragraphs once, collect the indexes of paragraphs that match
# either pattern, then delete them back-to-front (highest index first) so
# that removing one paragraph never shifts the index of another one still
# queued for deletion.

$d = $word.ActiveDocument

$toDelete = New-Object System.Collections.ArrayList

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $text = $rng.Text
    $len = $text.Length

    if ($len -lt 1) { continue }

    $lastCode = [int]$text[$len - 1]
    if ($lastCode -ne 13) { continue }   # paragraph must end with a real pilcrow (skip table row/cell-end markers)

    $before = $p.Format.SpaceBefore
    $after = $p.Format.SpaceAfter

    $isEmptySpacer = ($len -eq 1 -and $before -eq 2 -and $after -eq 4)

    $isSeparatorRule = $false
    if ($len -eq 61 -and $before -eq 6 -and $after -eq 6) {
        $isSeparatorRule = $true
        for ($j = 0; $j -lt ($len - 1); $j++) {
            if ([int]$text[$j] -ne 9472) { $isSeparatorRule = $false; break }
        }
    }

    if ($isEmptySpacer -or $isSeparatorRule) {
        [void]$toDelete.Add($i)
    }
}

for ($k = $toDelete.Count - 1; $k -ge 0; $k--) {
    $idx = $toDelete[$k]
    $d.Paragraphs.Item($idx).Range.Delete()
}

Write-Output "Removed $($toDelete.Count) paragraphs"
